# Modified correct coordinates for TU Delft Campus to Parking P1 (Aula).
# This updates the distance (km) and duration (minutes) matrices for the
# "TU Delft Campus" station/stop, which is stored in column N (rows 2-13)
# and mirrored in row 14 (columns B-M) of each matrix sheet.

$wb = $excel.ActiveWorkbook

# Sheet 0: Afstand_km (distances in km)
$ws0 = $wb.Worksheets.Item("Afstand_km")

# Sheet 1: Duur_minuten (durations in minutes)
$ws1 = $wb.Worksheets.Item("Duur_minuten")

# --- Afstand_km: column N (distance *to* TU Delft Campus) ---
$ws0.Range("N2").Value = 72.84
$ws0.Range("N3").Value = 12.69
$ws0.Range("N4").Value = 14.83
$ws0.Range("N5").Value = 66.01000000000001
$ws0.Range("N6").Value = 126.49
$ws0.Range("N7").Value = 94.58
$ws0.Range("N8").Value = 92.51000000000001
$ws0.Range("N9").Value = 239.63
$ws0.Range("N10").Value = 64.65000000000001
$ws0.Range("N11").Value = 131.73
$ws0.Range("N12").Value = 195.95
$ws0.Range("N13").Value = 233.15

# --- Afstand_km: row 14 (distance *from* TU Delft Campus) ---
$ws0.Range("B14").Value = 68.13
$ws0.Range("C14").Value = 12.4
$ws0.Range("D14").Value = 12.66
$ws0.Range("E14").Value = 65.33
$ws0.Range("F14").Value = 128.24
$ws0.Range("G14").Value = 96.09
$ws0.Range("H14").Value = 91.66
$ws0.Range("I14").Value = 239.5
$ws0.Range("J14").Value = 64.13
$ws0.Range("K14").Value = 131.37
$ws0.Range("L14").Value = 195.69
$ws0.Range("M14").Value = 233.41

# --- Duur_minuten: column N (duration *to* TU Delft Campus) ---
$ws1.Range("N2").Value = 3428.58
$ws1.Range("N3").Value = 1105.12
$ws1.Range("N4").Value = 931.77
$ws1.Range("N5").Value = 3119.69
$ws1.Range("N6").Value = 6065.56
$ws1.Range("N7").Value = 4599.27
$ws1.Range("N8").Value = 4193.88
$ws1.Range("N9").Value = 10156.5
$ws1.Range("N10").Value = 3536.67
$ws1.Range("N11").Value = 6262.38
$ws1.Range("N12").Value = 8523.120000000001
$ws1.Range("N13").Value = 9228.549999999999

# --- Duur_minuten: row 14 (duration *from* TU Delft Campus) ---
$ws1.Range("B14").Value = 3537.2
$ws1.Range("C14").Value = 1066.93
$ws1.Range("D14").Value = 1053.47
$ws1.Range("E14").Value = 3172.86
$ws1.Range("F14").Value = 6062.15
$ws1.Range("G14").Value = 4405.5
$ws1.Range("H14").Value = 4152.39
$ws1.Range("I14").Value = 10346.36
$ws1.Range("J14").Value = 3193.51
$ws1.Range("K14").Value = 5869.82
$ws1.Range("L14").Value = 8577.59
$ws1.Range("M14").Value = 9065.74
